$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.582.72"
$ws.Range("E2").Value = "  -3.74%  "
$ws.Range("D3").Value = "2.540.66"
$ws.Range("E3").Value = "  -3.53%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "507.81"
$ws.Range("E5").Value = "  -4.02%  "
$ws.Range("D6").Value = "143.96"
$ws.Range("E6").Value = "  -7.16%  "
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -4.35%  "
$ws.Range("D9").Value = "2.545.62"
$ws.Range("E9").Value = "  -3.67%  "
$ws.Range("D10").Value = "6.10"
$ws.Range("E10").Value = "  -8.58%  "
$ws.Range("E11").Value = "  -6.11%  "
$ws.Range("D12").Value = "0.332"
$ws.Range("E12").Value = "  -5.21%  "
$ws.Range("E13").Value = "  -0.56%  "
$ws.Range("D14").Value = "2.986.03"
$ws.Range("E14").Value = "  -3.54%  "
$ws.Range("D15").Value = "58.573.46"
$ws.Range("E15").Value = "  -3.78%  "
$ws.Range("D16").Value = "20.71"
$ws.Range("E16").Value = "  -5.83%  "
$ws.Range("E17").Value = "  -5.82%  "
$ws.Range("D18").Value = "2.544.55"
$ws.Range("E18").Value = "  -3.59%  "
$ws.Range("E19").Value = "  -4.77%  "
$ws.Range("D20").Value = "336.40"
$ws.Range("E20").Value = "  -4.71%  "
$ws.Range("D21").Value = "10.10"
$ws.Range("E21").Value = "  -4.82%  "
$ws.Range("E22").Value = "  -0.24%  "
$ws.Range("E23").Value = "  -4.35%  "
$ws.Range("D24").Value = "60.53"
$ws.Range("E24").Value = "  -1.70%  "
$ws.Range("D25").Value = "0.411"
$ws.Range("E25").Value = "  -4.51%  "
$ws.Range("D26").Value = "0.998"
$ws.Range("E27").Value = "  -4.64%  "
$ws.Range("D28").Value = "2.651.81"
$ws.Range("E28").Value = "  -3.57%  "
$ws.Range("D29").Value = "0.0₃0787"
$ws.Range("E29").Value = "  -8.91%  "
$ws.Range("E30").Value = "  -5.78%  "
$ws.Range("D32").Value = "149.83"
$ws.Range("E32").Value = "  -0.33%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "18.55"
$ws.Range("E33").Value = "  -4.71%  "
$ws.Range("B34").Value = "Aptos"
$ws.Range("C34").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D34").Value = "5.84"
$ws.Range("E34").Value = "  -5.07%  "
$ws.Range("E35").Value = "  -5.25%  "
$ws.Range("D36").Value = "0.919"
$ws.Range("E36").Value = "  +3.48%  "
$ws.Range("D37").Value = "3.90"
$ws.Range("E37").Value = "  -6.14%  "
$ws.Range("E38").Value = "  -7.19%  "
$ws.Range("D39").Value = "36.08"
$ws.Range("E39").Value = "  -1.33%  "
$ws.Range("D40").Value = "0.822"
$ws.Range("E40").Value = "  -11.08%  "
$ws.Range("E41").Value = "  -6.72%  "
$ws.Range("D42").Value = "283.78"
$ws.Range("E42").Value = "  -6.94%  "
$ws.Range("E43").Value = "  -7.77%  "
$ws.Range("D44").Value = "0.0995"
$ws.Range("E44").Value = "  -2.63%  "
$ws.Range("D45").Value = "0.998"
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").Value = "0.600"
$ws.Range("E46").Value = "  -6.29%  "
$ws.Range("E47").Value = "  -5.12%  "
$ws.Range("D48").Value = "18.69"
$ws.Range("E48").Value = "  -5.28%  "
$ws.Range("D49").Value = "10.30"
$ws.Range("E49").Value = "  -0.40%  "
$ws.Range("D50").Value = "0.0227"
$ws.Range("E50").Value = "  -4.82%  "
